# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Chirimoya" (Vega Modelo de Temuco) as
# row 84, pushing the existing rows 84-125 down to 85-126.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 84; this shifts rows 84..125
# down to 85..126 and extends the sheet's dimension accordingly.
$ws.Rows.Item(84).Insert()

# Populate the newly inserted row 84 with the new data point.
$ws.Cells.Item(84, 1).Value  = 10
$ws.Cells.Item(84, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(84, 3).Value  = "La Araucanía"
$ws.Cells.Item(84, 4).Value  = 44784
$ws.Cells.Item(84, 5).Value  = 9
$ws.Cells.Item(84, 6).Value  = "Fruta"
$ws.Cells.Item(84, 7).Value  = 100107
$ws.Cells.Item(84, 8).Value  = "Otros"
$ws.Cells.Item(84, 9).Value  = 100107002
$ws.Cells.Item(84, 10).Value = "Chirimoya"
$ws.Cells.Item(84, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(84, 12).Value = "Primera"
$ws.Cells.Item(84, 13).Value = 125
$ws.Cells.Item(84, 14).Value = 4000
$ws.Cells.Item(84, 15).Value = 4000
$ws.Cells.Item(84, 16).Value = 4000
$ws.Cells.Item(84, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(84, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(84, 19).Value = 4000
$ws.Cells.Item(84, 20).Value = 1
